# Insert 2 new rows at 516-517 (existing rows 516 onward shift down to 518 onward)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("516:517").Insert()

# New row 516 values
$ws.Range("A516").Value = 9
$ws.Range("B516").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C516").Value = "Metropolitana"
$ws.Range("D516").Value = 45021
$ws.Range("E516").Value = 13
$ws.Range("F516").Value = 100112052
$ws.Range("G516").Value = "Albahaca"
$ws.Range("H516").Value = "Sin especificar"
$ws.Range("I516").Value = "Primera"
$ws.Range("J516").Value = 260
$ws.Range("K516").Value = 3000
$ws.Range("L516").Value = 3000
$ws.Range("M516").Value = 3000
$ws.Range("N516").Value = "`$/docena de matas"
$ws.Range("O516").Value = "Región Metropolitana"
$ws.Range("P516").Value = 500
$ws.Range("Q516").Value = 6
$ws.Range("R516").Value = "Hortaliza"

# New row 517 values
$ws.Range("A517").Value = 9
$ws.Range("B517").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C517").Value = "Metropolitana"
$ws.Range("D517").Value = 45021
$ws.Range("E517").Value = 13
$ws.Range("F517").Value = 100112052
$ws.Range("G517").Value = "Albahaca"
$ws.Range("H517").Value = "Sin especificar"
$ws.Range("I517").Value = "Segunda"
$ws.Range("J517").Value = 170
$ws.Range("K517").Value = 2500
$ws.Range("L517").Value = 2500
$ws.Range("M517").Value = 2500
$ws.Range("N517").Value = "`$/docena de matas"
$ws.Range("O517").Value = "Región Metropolitana"
$ws.Range("P517").Value = 417
$ws.Range("Q517").Value = 6
$ws.Range("R517").Value = "Hortaliza"
